# Updated cryptos list on Fri Nov 29 07:49:10 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns, and a few reordered coin rows
# (Coin/Link/Price swapped between B/C/D for rows 37-38, 48-49, and 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price values must stay text (matching the source sheet's
# inline-string cells), so write them with a leading apostrophe (Excel's
# "treat as text" quote-prefix) and then strip the resulting cell style back
# to Normal so no stray number-format style lingers on the cell.

$ws.Range('D2').Value = '''95.853.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '''3.550.94'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''239.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').Value = '''650.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('D7').Value = '''1.61'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +9.96%  '
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('E9').Value = '  +6.19%  '
$ws.Range('D11').Value = '''3.548.79'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('D12').Value = '''43.09'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').Value = '''0.201'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '''4.209.41'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.54%  '
$ws.Range('D16').Value = '''95.756.91'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').Value = '''0.0000257'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '''3.537.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').Value = '''7.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.08%  '
$ws.Range('E20').Value = '  -2.34%  '
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = '''0.515'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.90%  '
$ws.Range('D23').Value = '''502.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('D24').Value = '''3.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.89%  '
$ws.Range('D25').Value = '''6.87'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.04%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '''95.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').Value = '''12.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('D29').Value = '''3.741.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').Value = '''0.152'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.89%  '
$ws.Range('E31').Value = '  -3.95%  '
$ws.Range('D32').Value = '''11.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').Value = '''31.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '''8.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.20%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''607.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.51%  '
$ws.Range('D39').Value = '''0.559'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').Value = '''1.61'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.28%  '
$ws.Range('D42').Value = '''0.150'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.21%  '
$ws.Range('D43').Value = '''0.896'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('E44').Value = '  +4.97%  '
$ws.Range('D45').Value = '''5.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').Value = '''23.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Value = '''2.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''33.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '''0.0417'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').Value = '''3.52'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '''52.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.87%  '
